# Updates the cryptos price/volume snapshot (GitHub Actions data refresh).
# D4 is a never-edited, unstyled "Price" cell - used only as a style donor so
# that forcing numeric-looking prices (e.g. "245.38") to stay text via
# NumberFormat "@" doesn't leave a stray cell style behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$defaultStyle = $ws.Range('D4').Style

$ws.Range('D2').Value = '35.285.26'
$ws.Range('E2').Value = '  +0.32%  '
$ws.Range('D3').Value = '1.881.41'
$ws.Range('E3').Value = '  -1.07%  '
$ws.Range('E4').Value = '  -0.61%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.38'
$ws.Range('D5').Style = $defaultStyle
$ws.Range('E5').Value = '  -3.04%  '
$ws.Range('E6').Value = '  -1.24%  '
$ws.Range('E7').Value = '  -0.68%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '43.40'
$ws.Range('D8').Style = $defaultStyle
$ws.Range('E8').Value = '  +5.71%  '
$ws.Range('E9').Value = '  -1.28%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '53.57'
$ws.Range('D10').Style = $defaultStyle
$ws.Range('E11').Value = '  -1.78%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0972'
$ws.Range('D12').Style = $defaultStyle
$ws.Range('E12').Value = '  -1.13%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '13.25'
$ws.Range('D13').Style = $defaultStyle
$ws.Range('E13').Value = '  +1.74%  '
$ws.Range('D14').Value = '2.156.11'
$ws.Range('E14').Value = '  -0.97%  '
$ws.Range('E15').Value = '  +2.69%  '
$ws.Range('E16').Value = '  -1.98%  '
$ws.Range('D17').Value = '1.876.32'
$ws.Range('E17').Value = '  -1.19%  '
$ws.Range('D18').Value = '35.413.96'
$ws.Range('E18').Value = '  +0.74%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '72.87'
$ws.Range('D19').Style = $defaultStyle
$ws.Range('E19').Value = '  -1.08%  '
$ws.Range('D20').Value = '0.0₃0819'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '244.17'
$ws.Range('D21').Style = $defaultStyle
$ws.Range('E21').Value = '  +0.55%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.76'
$ws.Range('D22').Style = $defaultStyle
$ws.Range('E22').Value = '  -1.62%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.93'
$ws.Range('D23').Style = $defaultStyle
$ws.Range('E23').Value = '  -2.53%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.65'
$ws.Range('D24').Style = $defaultStyle
$ws.Range('E24').Value = '  +9.39%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '165.58'
$ws.Range('D26').Style = $defaultStyle
$ws.Range('E26').Value = '  -0.85%  '
$ws.Range('B27').Value = 'PancakeSwap'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.12'
$ws.Range('D27').Style = $defaultStyle
$ws.Range('E27').Value = '  -6.74%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.24'
$ws.Range('D29').Style = $defaultStyle
$ws.Range('E29').Value = '  -1.42%  '
$ws.Range('E30').Value = '  -2.05%  '
$ws.Range('D31').Value = '4.128.46'
$ws.Range('E31').Value = '  +0.01%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.70'
$ws.Range('D32').Style = $defaultStyle
$ws.Range('E32').Value = '  +7.96%  '
$ws.Range('E33').Value = '  -1.44%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0584'
$ws.Range('D34').Style = $defaultStyle
$ws.Range('E34').Value = '  -4.84%  '
$ws.Range('E35').Value = '  -1.66%  '
$ws.Range('E36').Value = '  -0.67%  '
$ws.Range('E37').Value = '  -10.14%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.846'
$ws.Range('D38').Style = $defaultStyle
$ws.Range('E38').Value = '  -1.10%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.94'
$ws.Range('D39').Style = $defaultStyle
$ws.Range('E39').Value = '  -3.17%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0698'
$ws.Range('D40').Style = $defaultStyle
$ws.Range('E40').Value = '  +7.74%  '
$ws.Range('E41').Value = '  +2.07%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '17.21'
$ws.Range('D42').Style = $defaultStyle
$ws.Range('E42').Value = '  -1.05%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '96.32'
$ws.Range('D43').Style = $defaultStyle
$ws.Range('E43').Value = '  -6.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.07'
$ws.Range('D44').Style = $defaultStyle
$ws.Range('E44').Value = '  -2.71%  '
$ws.Range('D45').Value = '1.299.03'
$ws.Range('E45').Value = '  -1.78%  '
$ws.Range('E46').Value = '  -5.37%  '
$ws.Range('E47').Value = '  +7.12%  '
$ws.Range('E48').Value = '  -2.10%  '
$ws.Range('E49').Value = '  -0.83%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '12.20'
$ws.Range('D50').Style = $defaultStyle
$ws.Range('E50').Value = '  +2.21%  '
